$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Unprotect()

# Text fixes
$ws.Range("B51").Value = "D R Horton Inc"
$ws.Range("A59").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution." + [char]10 + "Model holdings provided as of 2021-05-17 for illustrative purposes only and are subject to change."

# Weight (D) and Percent Change (E) updates
$ws.Range("D2").Value = 0.01459890332567556
$ws.Range("E2").Value = -0.02084052964881955
$ws.Range("D3").Value = 0.0501618430332414
$ws.Range("E3").Value = 0.01473517639393096
$ws.Range("D4").Value = 0.01429010970323009
$ws.Range("E4").Value = -0.01339664974840438
$ws.Range("D5").Value = 0.009832835630854382
$ws.Range("E5").Value = -0.003039128783082101
$ws.Range("D6").Value = 0.01551457602707234
$ws.Range("E6").Value = -0.00286513404734301
$ws.Range("D7").Value = 0.02019133637282133
$ws.Range("E7").Value = -0.002028838489964491
$ws.Range("D8").Value = 0.004748666925169264
$ws.Range("E8").Value = 0.03280870004129777
$ws.Range("D9").Value = 0.006810019996939582
$ws.Range("E9").Value = 0.02336679282540732
$ws.Range("D10").Value = 0.01450252983221879
$ws.Range("E10").Value = -0.001255650426920996
$ws.Range("D11").Value = 0.008353614568494625
$ws.Range("E11").Value = -0.00409897153077976
$ws.Range("D12").Value = 0.01592777426742931
$ws.Range("E12").Value = 0.006097560975609762
$ws.Range("D13").Value = 0.002963422667017023
$ws.Range("E13").Value = 0.0117647058823529
$ws.Range("D14").Value = 0.005934004863586712
$ws.Range("E14").Value = 0.002885170225043421
$ws.Range("D15").Value = 0.01480556470263277
$ws.Range("E15").Value = 0.004024144869215096
$ws.Range("D16").Value = 0.01072435270210824
$ws.Range("E16").Value = 0.01123301985370939
$ws.Range("D17").Value = 0.02096683793693263
$ws.Range("E17").Value = 0.002137894174238264
$ws.Range("D18").Value = 0.008482361586885196
$ws.Range("E18").Value = 0.008939580764488086
$ws.Range("D19").Value = 0.01695575819763488
$ws.Range("E19").Value = 0.0009987075549289237
$ws.Range("D20").Value = 0.01221540255242458
$ws.Range("E20").Value = -0.003465674532388885
$ws.Range("D21").Value = 0.007426486619813585
$ws.Range("E21").Value = -0.02698511166253104
$ws.Range("D22").Value = 0.01470031961820856
$ws.Range("E22").Value = 0.002727380528874779
$ws.Range("D23").Value = 0.02004512632800013
$ws.Range("E23").Value = -0.003983228511530323
$ws.Range("D24").Value = 0.01029736458526504
$ws.Range("E24").Value = 0.01236083106864205
$ws.Range("D25").Value = 0.02016223132877027
$ws.Range("E25").Value = -0.002729608220937707
$ws.Range("D26").Value = 0.01414467786814287
$ws.Range("E26").Value = -0.0004159349292909287
$ws.Range("D27").Value = 0.020647118249822
$ws.Range("E27").Value = -0.01810975306389817
$ws.Range("D28").Value = 0.05554238513401758
$ws.Range("E28").Value = -0.009258532757944304
$ws.Range("D29").Value = 0.02081331272061195
$ws.Range("E29").Value = -0.002032520325203291
$ws.Range("D30").Value = 0.02905399349251059
$ws.Range("E30").Value = -0.01114253878460614
$ws.Range("D31").Value = 0.01499003153797671
$ws.Range("E31").Value = -0.002689204763734132
$ws.Range("D32").Value = 0.01315361220771758
$ws.Range("E32").Value = -0.001978417266187193
$ws.Range("D33").Value = 0.01786271495001102
$ws.Range("E33").Value = 0.008657465495608418
$ws.Range("D34").Value = 0.04255337984923264
$ws.Range("E34").Value = 0.004626093979055224
$ws.Range("D35").Value = 0.01092980007188043
$ws.Range("E35").Value = -0.003417634996582541
$ws.Range("D36").Value = 0.01014797944472918
$ws.Range("E36").Value = 0.003950871768444664
$ws.Range("D37").Value = 0.01081166783426142
$ws.Range("E37").Value = 0.009501187648456089
$ws.Range("D38").Value = 0.007502564403408014
$ws.Range("E38").Value = 0.003360716952949927
$ws.Range("D39").Value = 0.01226296673136578
$ws.Range("E39").Value = 0.008970727101038856
$ws.Range("D40").Value = 0.01745493304940258
$ws.Range("E40").Value = 0.0007561436672967048
$ws.Range("D41").Value = 0.01752553223646975
$ws.Range("E41").Value = -0.00492710583153344
$ws.Range("D42").Value = 0.03199282473193367
$ws.Range("E42").Value = -0.006947079599521233
$ws.Range("D43").Value = 0.01145393366623652
$ws.Range("E43").Value = -0.001869236091747384
$ws.Range("D44").Value = 0.02189925771180477
$ws.Range("E44").Value = -0.002203225522164454
$ws.Range("D45").Value = 0.01237222737801735
$ws.Range("E45").Value = -0.003804175535528254
$ws.Range("D46").Value = 0.008668851767537495
$ws.Range("E46").Value = -0.0005493972788674917
$ws.Range("D47").Value = 0.01357241355815924
$ws.Range("E47").Value = 0.01056846798497313
$ws.Range("D48").Value = 0.01084117754737416
$ws.Range("E48").Value = 0.006977285060857152
$ws.Range("D49").Value = 0.01600378979424502
$ws.Range("E49").Value = -0.00864775286799635
$ws.Range("D50").Value = 0.008694470431980303
$ws.Range("E50").Value = -0.004990870359099397
$ws.Range("D51").Value = 0.01138993992339142
$ws.Range("E51").Value = -0.00664658843078203
$ws.Range("D52").Value = 0.00829938891423051
$ws.Range("E52").Value = -0.005322221305388064
$ws.Range("D53").Value = 0.009865364797734978
$ws.Range("E53").Value = -0.005057979017117686
$ws.Range("D54").Value = 0.1352013031539116
$ws.Range("E54").Value = -0.0001970831690975006
$ws.Range("D55").Value = 0.04373694346945661
$ws.Range("E55").Value = -0.002903811252268795
$ws.Range("E56").Value = -0.0007901256814323299

$ws.Protect()
